# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary header values ---
# Valor Mora (total) and Cant. Trabajadores / Cant. Periodos counters
$ws.Range("E11").Value = 339606
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 10

# --- Replace the detail table (rows 16-29) with the refreshed dataset ---
$data = @(
    @("CC", "45646635",   "PAOLA MARGARITA RODRIGUEZ CABALLERO",  "1704", 7933,  850000),
    @("CC", "1044935054", "GUSTAVO ANDRES MARTINEZ PADILLA",      "2010", 39227, 980657),
    @("CC", "1044935054", "GUSTAVO ANDRES MARTINEZ PADILLA",      "2011", 39227, 980657),
    @("CC", "1044908417", "HAWYN JESSID CASTRO SARMIENTO",        "2012", 36612, 980657),
    @("CC", "1044912807", "MIRLEYDIS DE JESUS ROMERO HURTADO",    "2101", 7268,  908526),
    @("CC", "73353731",   "EULISE MATTAS BARRIOS",                "2105", 30284, 908526),
    @("CC", "3828458",    "LEIDER TOMAS LORA TOBIAS",              "2105", 29073, 908526),
    @("CC", "73353731",   "EULISE MATTAS BARRIOS",                "2106", 36341, 908526),
    @("CC", "3828458",    "LEIDER TOMAS LORA TOBIAS",              "2106", 36341, 908526),
    @("CC", "73353731",   "EULISE MATTAS BARRIOS",                "2107", 8480,  908526),
    @("CC", "3828458",    "LEIDER TOMAS LORA TOBIAS",              "2107", 8480,  908526),
    @("CC", "3828458",    "LEIDER TOMAS LORA TOBIAS",              "2109", 36341, 908526),
    @("CC", "1007323757", "CARLOS ALFREDO RODRIGUEZ CERVANTES",   "2204", 1333,  1000000),
    @("CC", "1051824856", "CARLOS ANDRES DE ORO BUSTILLO",        "2204", 22666, 1000000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}

# The last row of the table (now row 29, previously row 36) carries the
# heavier "closing" bottom border of the box. Re-apply that formatting
# (copied from the old closing row) onto the new last data row before the
# stale trailing rows are removed.
$ws.Range("B36:J36").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the now-stale trailing rows (old rows 30-36) so the
#     "firma" footer block shifts up from rows 41/42 to rows 34/35 ---
$ws.Range("B30:J36").EntireRow.Delete()
